$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 391, pushing existing rows 391-489 down to 393-491
$ws.Rows.Item(391).Resize(2, 1).EntireRow.Insert()

# Populate the newly inserted row 391 (Primera) with the new week's data
$ws.Range("A391").Value = 11
$ws.Range("B391").Value = "Vega Monumental Concepción"
$ws.Range("C391").Value = "Bíobío"
$ws.Range("D391").Value = 44798
$ws.Range("E391").Value = 8
$ws.Range("F391").Value = 100112020
$ws.Range("G391").Value = "Tomate"
$ws.Range("H391").Value = "Larga vida"
$ws.Range("I391").Value = "Primera"
$ws.Range("J391").Value = 600
$ws.Range("K391").Value = 8500
$ws.Range("L391").Value = 9000
$ws.Range("M391").Value = 8750
$ws.Range("N391").Value = "$/bandeja 18 kilos"
$ws.Range("O391").Value = "Región de Arica y Parinacota"
$ws.Range("P391").Value = 486
$ws.Range("Q391").Value = 18
$ws.Range("R391").Value = "Hortaliza"

# Populate the newly inserted row 392 (Segunda) with the new week's data
$ws.Range("A392").Value = 11
$ws.Range("B392").Value = "Vega Monumental Concepción"
$ws.Range("C392").Value = "Bíobío"
$ws.Range("D392").Value = 44798
$ws.Range("E392").Value = 8
$ws.Range("F392").Value = 100112020
$ws.Range("G392").Value = "Tomate"
$ws.Range("H392").Value = "Larga vida"
$ws.Range("I392").Value = "Segunda"
$ws.Range("J392").Value = 300
$ws.Range("K392").Value = 7500
$ws.Range("L392").Value = 7500
$ws.Range("M392").Value = 7500
$ws.Range("N392").Value = "$/bandeja 18 kilos"
$ws.Range("O392").Value = "Región de Arica y Parinacota"
$ws.Range("P392").Value = 417
$ws.Range("Q392").Value = 18
$ws.Range("R392").Value = "Hortaliza"
